# Generate Report for Handback
# Update the "Correspond Handoff Datetime" (column E) and
# "Correspond Handback DateTime" (column H) values for the first data
# row (the 4aac4d61-... file) on both the zh-cn and de-de report sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 00:46:56"
$wsZhCn.Range("H2").Value = "2016-03-23 00:47:19"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 00:47:00"
$wsDeDe.Range("H2").Value = "2016-03-23 00:47:25"
